# Femacal de La Calera - Acelga: weekly fruit/vegetable price update.
# A new observation is inserted at row 542 (pushing every existing row
# from 542 down through 657 to 543 through 658), and the new row 542 is
# populated with this week's reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 542, shifting rows 542:657 down to 543:658
# (styles/number formats travel with the shifted cells automatically).
$ws.Rows(542).Insert()

# Populate the newly inserted row 542 with the new data point.
$ws.Range("A542").Value = 3
$ws.Range("B542").Value = "Femacal de La Calera"
$ws.Range("C542").Value = "Coquimbo"
$ws.Range("D542").Value = 45258
$ws.Range("E542").Value = 5
$ws.Range("F542").Value = 100112009
$ws.Range("G542").Value = "Acelga"
$ws.Range("H542").Value = "Sin especificar"
$ws.Range("I542").Value = "Primera"
$ws.Range("J542").Value = 210
$ws.Range("K542").Value = 3500
$ws.Range("L542").Value = 4000
$ws.Range("M542").Value = 3762
$ws.Range("N542").Value = "$/docena de atados (6 kilos)"
$ws.Range("O542").Value = "Provincia de Quillota"
$ws.Range("P542").Value = 627
$ws.Range("Q542").Value = 6
$ws.Range("R542").Value = "Hortaliza"
